# Update README/Logboek with detailed security topics:
# - Fix "X BOLA" -> "X - BOLA" typo in the big security-topics note (cell E2)
# - Add two new rows detailing XSS and CSRF as separate topics to look into
# - Move selection/scroll position to reflect the newly added content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the small formatting typo ("X BOLA" -> "X - BOLA") inside the long
# note that lives in E2, keeping the rest of the text identical.
$notes = $ws.Range("E2").Value2
$notes = $notes -replace "X BOLA", "X - BOLA"
$ws.Range("E2").Value = $notes

# Add the two new rows that break out XSS and CSRF as their own topics.
$ws.Range("A10").Value = "XSS in detail"
$ws.Range("A11").Value = "CSRF in detail"

# Reflect the new selection/scroll state from the edit session.
[void]$ws.Range("A1").Select()
[void]$ws.Range("E11").Select()
